$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.398581624031067
$ws.Range("B1").Value = 1.803339242935181
$ws.Range("C1").Value = 1.506727695465088
$ws.Range("D1").Value = 2.246083736419678
$ws.Range("E1").Value = 3.384918689727783
